$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = 42602.010104166664
$ws.Range("B17").Value = "Named"
$ws.Range("C17").Value = 7631
$ws.Range("D17").Value = 4127
$ws.Range("E17").Value = 260
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 76
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0

# Row 18
$ws.Range("A18").Value = 42602.481921296298
$ws.Range("B18").Value = "Named"
$ws.Range("C18").Value = 9608
$ws.Range("D18").Value = 6400
$ws.Range("E18").Value = 430
$ws.Range("F18").Value = 121
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = 0
